$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header style (bold, bordered, centered) from an existing header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy date-column style (date/time number format) down column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A72").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match page margins used by the other sheets (0.75/0.75/1/1/0.5/0.5 in)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# --- Header row values ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$data = New-Object 'object[,]' 71,4
$data[0,0] = 44955.99999999999
$data[0,1] = 178
$data[0,2] = -108.4770210179462
$data[0,3] = 489.525339750763
$data[1,0] = 44983.99999999999
$data[1,1] = 180
$data[1,2] = -133.2880132060477
$data[1,3] = 493.2052808772409
$data[2,0] = 44990.99999999999
$data[2,1] = 181
$data[2,2] = -119.8383498556163
$data[2,3] = 480.8907480646865
$data[3,0] = 44997.99999999999
$data[3,1] = 181
$data[3,2] = -116.6149456920872
$data[3,3] = 496.5058194718281
$data[4,0] = 45004.99999999999
$data[4,1] = 182
$data[4,2] = -139.9653400831315
$data[4,3] = 480.1746361826896
$data[5,0] = 45011.99999999999
$data[5,1] = 182
$data[5,2] = -127.1288212177423
$data[5,3] = 493.9426436895843
$data[6,0] = 45018.99999999999
$data[6,1] = 183
$data[6,2] = -121.7821888395473
$data[6,3] = 494.0874194500787
$data[7,0] = 45025.99999999999
$data[7,1] = 183
$data[7,2] = -133.0881066550734
$data[7,3] = 504.5026057073436
$data[8,0] = 45046.99999999999
$data[8,1] = 185
$data[8,2] = -110.0473793110489
$data[8,3] = 487.8571471685806
$data[9,0] = 45053.99999999999
$data[9,1] = 185
$data[9,2] = -135.3379980224832
$data[9,3] = 517.4412720414932
$data[10,0] = 45060.99999999999
$data[10,1] = 186
$data[10,2] = -134.6061387098783
$data[10,3] = 501.6174697646909
$data[11,0] = 45074.99999999999
$data[11,1] = 187
$data[11,2] = -120.477522321279
$data[11,3] = 501.891276245012
$data[12,0] = 45081.99999999999
$data[12,1] = 187
$data[12,2] = -122.7535488500974
$data[12,3] = 493.7016795922974
$data[13,0] = 45088.99999999999
$data[13,1] = 188
$data[13,2] = -126.742141487532
$data[13,3] = 520.788723102531
$data[14,0] = 45095.99999999999
$data[14,1] = 188
$data[14,2] = -142.3161425807704
$data[14,3] = 497.2858357608378
$data[15,0] = 45102.99999999999
$data[15,1] = 188
$data[15,2] = -111.8452258440161
$data[15,3] = 509.9373584543036
$data[16,0] = 45109.99999999999
$data[16,1] = 189
$data[16,2] = -128.6722252965511
$data[16,3] = 502.1467759960259
$data[17,0] = 45130.99999999999
$data[17,1] = 190
$data[17,2] = -95.6655534283571
$data[17,3] = 510.4481273855375
$data[18,0] = 45137.99999999999
$data[18,1] = 191
$data[18,2] = -140.9200263369565
$data[18,3] = 505.3034557174548
$data[19,0] = 45151.99999999999
$data[19,1] = 192
$data[19,2] = -110.9773350679238
$data[19,3] = 502.4710672561189
$data[20,0] = 45158.99999999999
$data[20,1] = 192
$data[20,2] = -123.7865941739009
$data[20,3] = 497.3526526495955
$data[21,0] = 45165.99999999999
$data[21,1] = 193
$data[21,2] = -119.6528817817718
$data[21,3] = 489.8062248733799
$data[22,0] = 45172.99999999999
$data[22,1] = 193
$data[22,2] = -115.280201892116
$data[22,3] = 515.1156838472633
$data[23,0] = 45179.99999999999
$data[23,1] = 194
$data[23,2] = -117.1994121981795
$data[23,3] = 516.9035798958848
$data[24,0] = 45186.99999999999
$data[24,1] = 194
$data[24,2] = -115.605683895288
$data[24,3] = 513.6434648024842
$data[25,0] = 45193.99999999999
$data[25,1] = 195
$data[25,2] = -130.2717854501206
$data[25,3] = 504.6946858846759
$data[26,0] = 45200.99999999999
$data[26,1] = 195
$data[26,2] = -113.2246016814282
$data[26,3] = 502.2390206428643
$data[27,0] = 45221.99999999999
$data[27,1] = 197
$data[27,2] = -112.0215852827173
$data[27,3] = 507.6570966803857
$data[28,0] = 45228.99999999999
$data[28,1] = 197
$data[28,2] = -121.9915075529279
$data[28,3] = 500.2591337180123
$data[29,0] = 45235.99999999999
$data[29,1] = 198
$data[29,2] = -108.2123100283474
$data[29,3] = 512.755624095086
$data[30,0] = 45242.99999999999
$data[30,1] = 198
$data[30,2] = -118.3167515448005
$data[30,3] = 534.9696204861478
$data[31,0] = 45249.99999999999
$data[31,1] = 199
$data[31,2] = -101.7314647401948
$data[31,3] = 501.0361697318867
$data[32,0] = 45256.99999999999
$data[32,1] = 199
$data[32,2] = -109.3254627193016
$data[32,3] = 504.0468349202963
$data[33,0] = 45270.99999999999
$data[33,1] = 200
$data[33,2] = -78.56490822092861
$data[33,3] = 515.7816704654687
$data[34,0] = 45277.99999999999
$data[34,1] = 201
$data[34,2] = -102.7977106738188
$data[34,3] = 503.606232988604
$data[35,0] = 45298.99999999999
$data[35,1] = 202
$data[35,2] = -113.5116026231096
$data[35,3] = 523.0661654639345
$data[36,0] = 45312.99999999999
$data[36,1] = 203
$data[36,2] = -75.21943901497082
$data[36,3] = 519.6799673076135
$data[37,0] = 45326.99999999999
$data[37,1] = 204
$data[37,2] = -108.6629377152641
$data[37,3] = 534.9400180191346
$data[38,0] = 45375.99999999999
$data[38,1] = 207
$data[38,2] = -115.578448949801
$data[38,3] = 497.8124071465521
$data[39,0] = 45382.99999999999
$data[39,1] = 208
$data[39,2] = -87.84668340440138
$data[39,3] = 542.1010652774044
$data[40,0] = 45389.99999999999
$data[40,1] = 208
$data[40,2] = -94.82724891634869
$data[40,3] = 530.5855357389886
$data[41,0] = 45403.99999999999
$data[41,1] = 209
$data[41,2] = -97.74559513045141
$data[41,3] = 542.4732009537543
$data[42,0] = 45410.99999999999
$data[42,1] = 210
$data[42,2] = -84.04291551073086
$data[42,3] = 521.2361073910565
$data[43,0] = 45417.99999999999
$data[43,1] = 210
$data[43,2] = -86.39668430883771
$data[43,3] = 521.3424701078263
$data[44,0] = 45431.99999999999
$data[44,1] = 211
$data[44,2] = -126.0645758816462
$data[44,3] = 515.8335238244993
$data[45,0] = 45438.99999999999
$data[45,1] = 212
$data[45,2] = -124.7440305714247
$data[45,3] = 509.3121640557217
$data[46,0] = 45445.99999999999
$data[46,1] = 212
$data[46,2] = -113.9730260655633
$data[46,3] = 517.1897348388418
$data[47,0] = 45459.99999999999
$data[47,1] = 213
$data[47,2] = -123.6734187896474
$data[47,3] = 514.1018811837342
$data[48,0] = 45466.99999999999
$data[48,1] = 214
$data[48,2] = -120.0869825673987
$data[48,3] = 516.294941904277
$data[49,0] = 45487.99999999999
$data[49,1] = 215
$data[49,2] = -147.1312701001027
$data[49,3] = 527.0378143217013
$data[50,0] = 45501.99999999999
$data[50,1] = 216
$data[50,2] = -114.4353025055987
$data[50,3] = 528.1442825207286
$data[51,0] = 45515.99999999999
$data[51,1] = 217
$data[51,2] = -102.5255024035436
$data[51,3] = 533.4184505771333
$data[52,0] = 45522.99999999999
$data[52,1] = 217
$data[52,2] = -101.9834459738002
$data[52,3] = 514.5418058107175
$data[53,0] = 45529.99999999999
$data[53,1] = 218
$data[53,2] = -91.75892864852176
$data[53,3] = 530.199081049958
$data[54,0] = 45536.99999999999
$data[54,1] = 218
$data[54,2] = -86.79608211444467
$data[54,3] = 542.5092961691141
$data[55,0] = 45543.99999999999
$data[55,1] = 219
$data[55,2] = -87.63691037322263
$data[55,3] = 521.3423802117663
$data[56,0] = 45550.99999999999
$data[56,1] = 219
$data[56,2] = -81.76853774964029
$data[56,3] = 545.1258627726694
$data[57,0] = 45557.99999999999
$data[57,1] = 220
$data[57,2] = -77.91751333647079
$data[57,3] = 514.2926960842901
$data[58,0] = 45564.99999999999
$data[58,1] = 220
$data[58,2] = -62.91798615888551
$data[58,3] = 528.9290973594202
$data[59,0] = 45578.99999999999
$data[59,1] = 221
$data[59,2] = -73.82393761985911
$data[59,3] = 535.0969738715339
$data[60,0] = 45585.99999999999
$data[60,1] = 222
$data[60,2] = -94.592459399516
$data[60,3] = 552.0105460764755
$data[61,0] = 45592.99999999999
$data[61,1] = 222
$data[61,2] = -66.38248415533093
$data[61,3] = 531.7058806335149
$data[62,0] = 45599.99999999999
$data[62,1] = 223
$data[62,2] = -90.51872123868819
$data[62,3] = 547.3616328795513
$data[63,0] = 45606.99999999999
$data[63,1] = 223
$data[63,2] = -102.5225905975624
$data[63,3] = 541.9913253504294
$data[64,0] = 45613.99999999999
$data[64,1] = 224
$data[64,2] = -70.67995274417912
$data[64,3] = 548.0923776262352
$data[65,0] = 45620.99999999999
$data[65,1] = 224
$data[65,2] = -76.35618848703504
$data[65,3] = 535.0116293634791
$data[66,0] = 45627.99999999999
$data[66,1] = 225
$data[66,2] = -78.82814403377132
$data[66,3] = 540.5579942679954
$data[67,0] = 45634.99999999999
$data[67,1] = 225
$data[67,2] = -95.66106338438604
$data[67,3] = 531.7250634175871
$data[68,0] = 45641.99999999999
$data[68,1] = 226
$data[68,2] = -88.52153779167267
$data[68,3] = 536.2020017355978
$data[69,0] = 45648.99999999999
$data[69,1] = 226
$data[69,2] = -77.10849518604748
$data[69,3] = 544.1805216785883
$data[70,0] = 45655.99999999999
$data[70,1] = 227
$data[70,2] = -86.26799267188312
$data[70,3] = 539.5554830332453

$wsForecast.Range("A2:D72").Value = $data

$wsForecast.Range("A1").Select()